# Append: 2025-11-16 12:31 JST
# - Bump the "取得日時" (fetched-at) timestamp on every existing row.
# - Refresh row 3 with a new job listing (title/price/url/score/skills).
# - Append 5 new job listings (rows 5-9), pushing the old row 5
#   ("paperspace...") down to row 10.
# - Widen column B a bit to fit the new titles.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-11-16 12:31:55"

# ---------------------------------------------------------------------
# 1) Move the old row 5 ("paperspace...") down to row 10 before we
#    overwrite row 5 with new content. (Values are hard-coded here,
#    rather than copied from A5:G5 at runtime, because this host's
#    Range.Value getter doesn't reliably round-trip through COM.)
# ---------------------------------------------------------------------
$ws.Range("B10").Value = "paperspaceで、comfyuiが動くようにして欲しい。"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5434935"
$ws.Range("G10").Value = 10

# ---------------------------------------------------------------------
# 2) Bump the timestamp on every row (2-10) to the new fetch time.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A8").Value = $newTimestamp
$ws.Range("A9").Value = $newTimestamp
$ws.Range("A10").Value = $newTimestamp

# ---------------------------------------------------------------------
# 3) Row 3 gets replaced with a different job listing.
# ---------------------------------------------------------------------
$ws.Range("B3").Value = "【在宅】英語ニュース → 日本語AI要約のオペレーション作業"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5435071"
$ws.Range("G3").Value = 295
$ws.Range("H3").Value = "🔥AI,Ai"

# ---------------------------------------------------------------------
# 4) New row 5 (replacing the old paperspace row which is now at 10).
# ---------------------------------------------------------------------
$ws.Range("B5").Value = "完全在宅GASエンジニア募集/課題テストからご依頼/時給1,163円~業務フロー効率化をお任せします"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5435126"
$ws.Range("G5").Value = 70
$ws.Range("H5").Value = "◆効率化"

# ---------------------------------------------------------------------
# 5) New row 6.
# ---------------------------------------------------------------------
$ws.Range("B6").Value = "【動画解析】スポーツ動作の“微妙な違い”を検出・可視化する仕組みの開発者を募集します"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5422314"
$ws.Range("G6").Value = 68
$ws.Range("H6").Value = "◆開発"

# ---------------------------------------------------------------------
# 6) New row 7.
# ---------------------------------------------------------------------
$ws.Range("B7").Value = "【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5435079"
$ws.Range("G7").Value = 25

# ---------------------------------------------------------------------
# 7) New row 8.
# ---------------------------------------------------------------------
$ws.Range("B8").Value = "【リーダー募集×リモートOK】SRE/インフラエンジニア(Google Cloud/長期金融系案件)"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5435080"
$ws.Range("G8").Value = 25

# ---------------------------------------------------------------------
# 8) New row 9.
# ---------------------------------------------------------------------
$ws.Range("B9").Value = "【ITエンジニア必見】セキュリティ基礎研修講師募集|報酬8万円"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5435101"
$ws.Range("G9").Value = 18

# ---------------------------------------------------------------------
# 9) Rebuild every hyperlink in column F (F2:F10) from scratch. Deleting
#    via a single range's Hyperlinks collection clears the whole sheet's
#    hyperlink set in this host, so re-add all of them in ref order --
#    that reproduces rId1..rId9 matching F2..F10 in column order, and
#    lets the changed F3 URL land on a clean, single relationship.
# ---------------------------------------------------------------------
$ws.Range("F3").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5434943")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5435071")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5434860")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5435126")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5422314")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5435079")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5435080")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5435101")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5434935")

$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("F10").Style = "Hyperlink"

# ---------------------------------------------------------------------
# 10) Widen column B (47 -> 52 chars). This host stores ColumnWidth with
#     a +0.8333 padding baked in, so feed it a value that lands exactly
#     on 52 once serialised.
# ---------------------------------------------------------------------
$ws.Range("B:B").ColumnWidth = 51.2
